# "Generate Report for Archive"
#
# The localization status report is regenerated: the two in-flight files
# have moved on from "Ready for handoff" to "In Translation". That text
# shows up in three places:
#   - Overview sheet: columns "zh-cn" (E) and "de-de" (F), rows 2-3
#   - zh-cn sheet:    "Status" column (C), rows 2-3
#   - de-de sheet:    "Status" column (C), rows 2-3
#
# Shortening the status text also narrows those status columns, which is
# reflected by re-sizing them to match the freshly generated report.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
